# Fruta / hortaliza, semanal
# Reassign the per-row record data (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg o Unidades) among rows 2-30: each destination row
# ends up with the values that previously belonged to a different row
# (a row-shuffle), while Mercado/Region/Fecha-codes/Categoria/Variedad and
# Clasificacion columns (A,B,C,E,F,G,H,R) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values travel together with a "record" when rows are reshuffled.
$cols = @("D", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Mapping: destination row -> source row (source row's old values are copied
# into the destination row).
$rowMap = [ordered]@{
    2  = 2
    3  = 12
    4  = 3
    5  = 4
    6  = 27
    7  = 21
    8  = 24
    9  = 19
    10 = 11
    11 = 20
    12 = 5
    13 = 6
    14 = 30
    15 = 15
    16 = 23
    17 = 9
    18 = 26
    19 = 22
    20 = 16
    21 = 29
    22 = 13
    23 = 28
    24 = 17
    25 = 25
    26 = 7
    27 = 14
    28 = 10
    29 = 18
    30 = 8
}

# Snapshot all current values first, since sources and destinations overlap.
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Now write the shuffled values back.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
